$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "External Data Source" (col D) and updated "Main Subject"
# (col C) entries for each notebook row. Order chosen to match the
# original authoring sequence (shared-string append order).
$ws.Range("C2").Value = "DS Industry / Country"
$ws.Range("C3").Value = "Men / Women"
$ws.Range("D2").Value = "Job, Funding, Global Data"
$ws.Range("D3").Value = "None?"
$ws.Range("D4").Value = "Job, Blog"
$ws.Range("C4").Value = "Data Scientist / Analyst"
$ws.Range("C5").Value = "Early Career Kaggler"
$ws.Range("C6").Value = "Cloud Computing"
$ws.Range("D6").Value = "Articles"
$ws.Range("D5").Value = "None?"

# Update column widths (values chosen so the stored/quantized width lands
# as close as possible to the target 4.875 / 44.75 / 28.625 / 38.625)
$ws.Columns.Item(1).ColumnWidth = 4.142857142857143
$ws.Columns.Item(2).ColumnWidth = 44.0
$ws.Columns.Item(3).ColumnWidth = 27.857142857142858
$ws.Columns.Item(4).ColumnWidth = 37.857142857142854

# Update selection to match author's final cursor position
$ws.Range("D10").Select()
